# Apply the "Elimna EC anteriores y se agregan nuevos, se modifica base de datos" edit:
# The period/mora rows (16-23) are reversed in order (periods 1910..2005 become 2005..1910),
# keeping each period's "Valor Mora" paired value attached to it. Since most Valor Mora
# values are identical (48000) except for the first/last rows (14400/30400), the net
# observable effect is that column E (Periodo Mora) and column F (Valor Mora) for rows
# 16-23 get reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Capture current values for the Periodo Mora (E) and Valor Mora (F) columns, rows 16-23
# (use Value2 -- Value returns an unusable placeholder in this environment)
$rows = 16..23
$periodos = @()
$valores = @()
foreach ($r in $rows) {
    $periodos += $ws.Range("E$r").Value2
    $valores  += $ws.Range("F$r").Value2
}

# Write them back in reverse order
$n = $rows.Count
for ($i = 0; $i -lt $n; $i++) {
    $r = $rows[$i]
    $ws.Range("E$r").Value2 = $periodos[$n - 1 - $i]
    $ws.Range("F$r").Value2 = $valores[$n - 1 - $i]
}
